$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Akash" row with allocation values
$ws.Range("B4").Value = 100
$ws.Range("D4").Value = 0

# Update the active selection to C4 as seen after the edit
$ws.Range("C4").Select()
